$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (A1): text change "tracado_via" -> "Cruzamento", bold formatting removed
$ws.Range("A1").Value = "Cruzamento"
$ws.Range("A1").Font.Bold = $false

# Insert a new row at position 5 (pushes "Não Informado"... down one row) and
# populate it with the new "n/a" value
$ws.Rows(5).Insert()
$ws.Range("A5").Value = "n/a"

# Restore the selection to the newly inserted row, matching post-insert UX
$null = $ws.Rows(5).Select()
